# Generate Report for Handback
# Fills in the handoff/handback timestamps for the 3de4c1a0 file once its
# xliff round-trip completed, and refreshes the de-de handoff timestamp
# for the 19201f17 file (report regeneration run).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 3de4c1a0-...md
$overview.Range("G3").Value = "2016-09-03 02:52:40"

# zh-cn sheet, row for 3de4c1a0-...md: Correspond Handoff Datetime / Correspond Handback DateTime
$zhcn.Range("H3").Value = "2016-09-03 02:52:36"
$zhcn.Range("K3").Value = "2016-09-03 02:52:52"

# de-de sheet, row for 19201f17-...md stays the same; row for 3de4c1a0-...md gets new
# Correspond Handoff Datetime / Correspond Handback DateTime
$dede.Range("H3").Value = "2016-09-03 02:52:40"
$dede.Range("K3").Value = "2016-09-03 02:52:59"
